$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1195.6666 # H28: was 1708.5
$ws.Cells.Item(28, 9).Value = 444 # I28: was 626.6667
$ws.Cells.Item(28, 11).Value = 444 # K28: was 626.6667
$ws.Cells.Item(28, 13).Value = 41 # M28: was -141.6667

$ws.Cells.Item(57, 8).Value = 44597.6 # H57: was 50191.5
$ws.Cells.Item(57, 10).Value = 44597.6 # J57: was 50191.5
$ws.Cells.Item(57, 12).Value = 133792.8 # L57: was 150574.5
$ws.Cells.Item(57, 14).Value = -134790.8 # N57: was -151572.5

$ws.Cells.Item(112, 8).Value = 1289.6666 # H112: was 1296.7894
$ws.Cells.Item(112, 10).Value = 1289.7106 # J112: was 1297.027
$ws.Cells.Item(112, 12).Value = 3869.1318 # L112: was 3891.081
$ws.Cells.Item(112, 14).Value = -6085.1318 # N112: was -6107.081

$ws.Cells.Item(113, 8).Value = 128397.25 # H113: was 334791.34
$ws.Cells.Item(113, 9).Value = 334669.66 # I113: was 501002.5
$ws.Cells.Item(113, 10).Value = 4633.8 # J113: was 2369
$ws.Cells.Item(113, 11).Value = 334669.66 # K113: was 501002.5
$ws.Cells.Item(113, 12).Value = 4633.8 # L113: was 2369
$ws.Cells.Item(113, 13).Value = -331415.66 # M113: was -497748.5
$ws.Cells.Item(113, 14).Value = -11141.8 # N113: was -8877

$ws.Cells.Item(132, 8).Value = 1821.262 # H132: was 1833.3572
$ws.Cells.Item(132, 9).Value = 1497.421 # I132: was 1510.7894
$ws.Cells.Item(132, 11).Value = 4492.263 # K132: was 4532.3682
$ws.Cells.Item(132, 13).Value = -1962.263 # M132: was -2002.3682

$ws.Cells.Item(137, 8).Value = 2056.6316 # H137: was 1911.1364
$ws.Cells.Item(137, 9).Value = 2092.125 # I137: was 1969.7222
$ws.Cells.Item(137, 10).Value = 1867.3334 # J137: was 1647.5
$ws.Cells.Item(137, 11).Value = 6276.375 # K137: was 5909.1666
$ws.Cells.Item(137, 12).Value = 5602.0002 # L137: was 4942.5
$ws.Cells.Item(137, 13).Value = -3726.375 # M137: was -3359.1666
$ws.Cells.Item(137, 14).Value = -10702.0002 # N137: was -10042.5

$ws.Cells.Item(138, 8).Value = 1670.6333 # H138: was 1710.5172
$ws.Cells.Item(138, 9).Value = 1049.8636 # I138: was 1075.381
$ws.Cells.Item(138, 11).Value = 3149.5908 # K138: was 3226.143
$ws.Cells.Item(138, 13).Value = 1990.4092 # M138: was 1913.857

$ws.Cells.Item(141, 8).Value = 9359.4 # H141: was 6281.846
$ws.Cells.Item(141, 9).Value = 9399.333000000001 # I141: was 5733.1816
$ws.Cells.Item(141, 11).Value = 28197.999 # K141: was 17199.5448
$ws.Cells.Item(141, 13).Value = -23017.999 # M141: was -12019.5448

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2319.375 # H2: was 2325
$ws.Cells.Item(2, 9).Value = 2074.3333 # I2: was 2080.3333
$ws.Cells.Item(2, 11).Value = 2074.3333 # K2: was 2080.3333
$ws.Cells.Item(2, 13).Value = -1961.3333 # M2: was -1967.3333

$ws.Cells.Item(45, 8).Value = 5244.8286 # H45: was 5360.8823
$ws.Cells.Item(45, 9).Value = 6491.905 # I45: was 6751.55
$ws.Cells.Item(45, 11).Value = 6491.905 # K45: was 6751.55
$ws.Cells.Item(45, 13).Value = -6114.905 # M45: was -6374.55

$ws.Cells.Item(102, 8).Value = 2202.0417 # H102: was 2393.8462
$ws.Cells.Item(102, 9).Value = 1168.6666 # I102: was 1156.091
$ws.Cells.Item(102, 10).Value = 9435.666999999999 # J102: was 9201.5
$ws.Cells.Item(102, 11).Value = 1168.6666 # K102: was 1156.091
$ws.Cells.Item(102, 12).Value = 9435.666999999999 # L102: was 9201.5
$ws.Cells.Item(102, 13).Value = 453.3334 # M102: was 465.9090000000001
$ws.Cells.Item(102, 14).Value = -12679.667 # N102: was -12445.5

$ws.Cells.Item(116, 8).Value = 2319.375 # H116: was 2325
$ws.Cells.Item(116, 9).Value = 2074.3333 # I116: was 2080.3333
$ws.Cells.Item(116, 11).Value = 2074.3333 # K116: was 2080.3333
$ws.Cells.Item(116, 13).Value = 219.6667000000002 # M116: was 213.6667000000002

$ws.Cells.Item(122, 8).Value = 2010.3438 # H122: was 1991.3871
$ws.Cells.Item(122, 9).Value = 2017.1613 # I122: was 1997.8
$ws.Cells.Item(122, 11).Value = 6051.4839 # K122: was 5993.4
$ws.Cells.Item(122, 13).Value = -3601.4839 # M122: was -3543.4

$ws.Cells.Item(132, 8).Value = 1170.6666 # H132: was 1191.4529
$ws.Cells.Item(132, 9).Value = 1136.96 # I132: was 1158.7551
$ws.Cells.Item(132, 11).Value = 3410.88 # K132: was 3476.2653
$ws.Cells.Item(132, 13).Value = -880.8800000000001 # M132: was -946.2653

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2319.375 # H3: was 2325
$ws.Cells.Item(3, 9).Value = 2074.3333 # I3: was 2080.3333
$ws.Cells.Item(3, 11).Value = 2074.3333 # K3: was 2080.3333
$ws.Cells.Item(3, 13).Value = -1960.3333 # M3: was -1966.3333

$ws.Cells.Item(20, 8).Value = 4996.853 # H20: was 5067.5
$ws.Cells.Item(20, 10).Value = 5851.857 # J20: was 6023.4287
$ws.Cells.Item(20, 12).Value = 5851.857 # L20: was 6023.4287
$ws.Cells.Item(20, 14).Value = -6345.857 # N20: was -6517.4287

$ws.Cells.Item(86, 8).Value = 2052.818 # H86: was 1493.138
$ws.Cells.Item(86, 9).Value = 1919.7142 # I86: was 1381.5834
$ws.Cells.Item(86, 10).Value = 2285.75 # J86: was 2028.6
$ws.Cells.Item(86, 11).Value = 1919.7142 # K86: was 1381.5834
$ws.Cells.Item(86, 12).Value = 2285.75 # L86: was 2028.6
$ws.Cells.Item(86, 13).Value = -796.7141999999999 # M86: was -258.5834
$ws.Cells.Item(86, 14).Value = -4531.75 # N86: was -4274.6

$ws.Cells.Item(89, 8).Value = 2052.818 # H89: was 1493.138
$ws.Cells.Item(89, 9).Value = 1919.7142 # I89: was 1381.5834
$ws.Cells.Item(89, 10).Value = 2285.75 # J89: was 2028.6
$ws.Cells.Item(89, 11).Value = 9598.571 # K89: was 6907.916999999999
$ws.Cells.Item(89, 12).Value = 11428.75 # L89: was 10143
$ws.Cells.Item(89, 13).Value = -3982.571 # M89: was -1291.916999999999
$ws.Cells.Item(89, 14).Value = -22660.75 # N89: was -21375

$ws.Cells.Item(141, 8).Value = 87769.39999999999 # H141: was 87771.5
$ws.Cells.Item(141, 10).Value = 87769.39999999999 # J141: was 87771.5
$ws.Cells.Item(141, 12).Value = 87769.39999999999 # L141: was 87771.5
$ws.Cells.Item(141, 14).Value = -98129.39999999999 # N141: was -98131.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1849.4 # H16: was 1049.3077
$ws.Cells.Item(16, 9).Value = 1849.4 # I16: was 1176.2
$ws.Cells.Item(16, 10).Value = 0 # J16: was 626.3333
$ws.Cells.Item(16, 11).Value = 1849.4 # K16: was 1176.2
$ws.Cells.Item(16, 12).Value = 0 # L16: was 626.3333
$ws.Cells.Item(16, 13).Value = -1562.4 # M16: was -889.2
$ws.Cells.Item(16, 14).ClearContents() # N16: was -1200.3333

$ws.Cells.Item(31, 8).Value = 13530.647 # H31: was 13014.973
$ws.Cells.Item(31, 10).Value = 24503 # J31: was 22252.5
$ws.Cells.Item(31, 12).Value = 24503 # L31: was 22252.5
$ws.Cells.Item(31, 14).Value = -25093 # N31: was -22842.5

$ws.Cells.Item(34, 8).Value = 13530.647 # H34: was 13014.973
$ws.Cells.Item(34, 10).Value = 24503 # J34: was 22252.5
$ws.Cells.Item(34, 12).Value = 24503 # L34: was 22252.5
$ws.Cells.Item(34, 14).Value = -24907 # N34: was -22656.5

$ws.Cells.Item(113, 8).Value = 1849.4 # H113: was 1049.3077
$ws.Cells.Item(113, 9).Value = 1849.4 # I113: was 1176.2
$ws.Cells.Item(113, 10).Value = 0 # J113: was 626.3333
$ws.Cells.Item(113, 11).Value = 1849.4 # K113: was 1176.2
$ws.Cells.Item(113, 12).Value = 0 # L113: was 626.3333
$ws.Cells.Item(113, 13).Value = 320.5999999999999 # M113: was 993.8
$ws.Cells.Item(113, 14).ClearContents() # N113: was -4966.3333

$ws.Cells.Item(141, 8).Value = 355346 # H141: was 410815.2
$ws.Cells.Item(141, 10).Value = 355346 # J141: was 410815.2
$ws.Cells.Item(141, 12).Value = 355346 # L141: was 410815.2
$ws.Cells.Item(141, 14).Value = -365706 # N141: was -421175.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 9333 # H62: was 340999.34
$ws.Cells.Item(62, 9).Value = 5000 # I62: was 999999
$ws.Cells.Item(62, 11).Value = 15000 # K62: was 2999997
$ws.Cells.Item(62, 13).Value = -14314 # M62: was -2999311

$ws.Cells.Item(64, 8).Value = 2323 # H64: was 2456.8572
$ws.Cells.Item(64, 9).Value = 687.5 # I64: was 766.6667
$ws.Cells.Item(64, 10).Value = 2917.7273 # J64: was 2917.818
$ws.Cells.Item(64, 11).Value = 2062.5 # K64: was 2300.0001
$ws.Cells.Item(64, 12).Value = 8753.1819 # L64: was 8753.454000000002
$ws.Cells.Item(64, 13).Value = -1792.5 # M64: was -2030.0001
$ws.Cells.Item(64, 14).Value = -9293.1819 # N64: was -9293.454000000002

$ws.Cells.Item(65, 8).Value = 9333 # H65: was 340999.34
$ws.Cells.Item(65, 9).Value = 5000 # I65: was 999999
$ws.Cells.Item(65, 11).Value = 45000 # K65: was 8999991
$ws.Cells.Item(65, 13).Value = -41568 # M65: was -8996559

$ws.Cells.Item(67, 8).Value = 2323 # H67: was 2456.8572
$ws.Cells.Item(67, 9).Value = 687.5 # I67: was 766.6667
$ws.Cells.Item(67, 10).Value = 2917.7273 # J67: was 2917.818
$ws.Cells.Item(67, 11).Value = 2062.5 # K67: was 2300.0001
$ws.Cells.Item(67, 12).Value = 8753.1819 # L67: was 8753.454000000002
$ws.Cells.Item(67, 13).Value = -1126.5 # M67: was -1364.0001
$ws.Cells.Item(67, 14).Value = -10625.1819 # N67: was -10625.454

$ws.Cells.Item(80, 8).Value = 1799.5 # H80: was 1800
$ws.Cells.Item(80, 10).Value = 1799.5 # J80: was 1800
$ws.Cells.Item(80, 12).Value = 5398.5 # L80: was 5400
$ws.Cells.Item(80, 14).Value = -7270.5 # N80: was -7272

$ws.Cells.Item(83, 8).Value = 1799.5 # H83: was 1800
$ws.Cells.Item(83, 10).Value = 1799.5 # J83: was 1800
$ws.Cells.Item(83, 12).Value = 16195.5 # L83: was 16200
$ws.Cells.Item(83, 14).Value = -25555.5 # N83: was -25560

$ws.Cells.Item(88, 8).Value = 10332.667 # H88: was 9249.75
$ws.Cells.Item(88, 10).Value = 10332.667 # J88: was 9249.75
$ws.Cells.Item(88, 12).Value = 30998.001 # L88: was 27749.25
$ws.Cells.Item(88, 14).Value = -31854.001 # N88: was -28605.25

$ws.Cells.Item(91, 8).Value = 10332.667 # H91: was 9249.75
$ws.Cells.Item(91, 10).Value = 10332.667 # J91: was 9249.75
$ws.Cells.Item(91, 12).Value = 30998.001 # L91: was 27749.25
$ws.Cells.Item(91, 14).Value = -33962.001 # N91: was -30713.25

$ws.Cells.Item(104, 8).Value = 1370.7142 # H104: was 1180.25
$ws.Cells.Item(104, 10).Value = 1174.75 # J104: was 909.2
$ws.Cells.Item(104, 12).Value = 3524.25 # L104: was 2727.6
$ws.Cells.Item(104, 14).Value = -8766.25 # N104: was -7969.6

$ws.Cells.Item(105, 8).Value = 15000 # H105: was 0
$ws.Cells.Item(105, 10).Value = 15000 # J105: was 0
$ws.Cells.Item(105, 12).Value = 45000 # L105: was 0
$ws.Cells.Item(105, 14).Value = -50242 # N105: was None

$ws.Cells.Item(113, 8).Value = 765.70966 # H113: was 760.03125
$ws.Cells.Item(113, 10).Value = 815.3684 # J113: was 803.8
$ws.Cells.Item(113, 12).Value = 2446.1052 # L113: was 2411.4
$ws.Cells.Item(113, 14).Value = -6786.1052 # N113: was -6751.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4719.2 # H80: was 4487.5884
$ws.Cells.Item(80, 9).Value = 3229.4 # I80: was 3229.8
$ws.Cells.Item(80, 10).Value = 5464.1 # J80: was 5011.6665
$ws.Cells.Item(80, 11).Value = 3229.4 # K80: was 3229.8
$ws.Cells.Item(80, 12).Value = 5464.1 # L80: was 5011.6665
$ws.Cells.Item(80, 13).Value = -2231.4 # M80: was -2231.8
$ws.Cells.Item(80, 14).Value = -7460.1 # N80: was -7007.6665

$ws.Cells.Item(83, 8).Value = 4719.2 # H83: was 4487.5884
$ws.Cells.Item(83, 9).Value = 3229.4 # I83: was 3229.8
$ws.Cells.Item(83, 10).Value = 5464.1 # J83: was 5011.6665
$ws.Cells.Item(83, 11).Value = 16147 # K83: was 16149
$ws.Cells.Item(83, 12).Value = 27320.5 # L83: was 25058.3325
$ws.Cells.Item(83, 13).Value = -11155 # M83: was -11157
$ws.Cells.Item(83, 14).Value = -37304.5 # N83: was -35042.3325

$ws.Cells.Item(102, 8).Value = 1861.5 # H102: was 1496.5
$ws.Cells.Item(102, 9).Value = 1913.4706 # I102: was 1519.0435
$ws.Cells.Item(102, 11).Value = 1913.4706 # K102: was 1519.0435
$ws.Cells.Item(102, 13).Value = -291.4706000000001 # M102: was 102.9565

$ws.Cells.Item(122, 8).Value = 1472.6666 # H122: was 1464.8
$ws.Cells.Item(122, 9).Value = 1472.6666 # I122: was 1531
$ws.Cells.Item(122, 10).Value = 0 # J122: was 1200
$ws.Cells.Item(122, 11).Value = 4417.9998 # K122: was 4593
$ws.Cells.Item(122, 12).Value = 0 # L122: was 3600
$ws.Cells.Item(122, 13).Value = -1967.9998 # M122: was -2143
$ws.Cells.Item(122, 14).ClearContents() # N122: was -8500

$ws.Cells.Item(126, 8).Value = 3466.0833 # H126: was 3463.6428
$ws.Cells.Item(126, 9).Value = 2832.5 # I126: was 2856.2856
$ws.Cells.Item(126, 10).Value = 4099.6665 # J126: was 4071
$ws.Cells.Item(126, 11).Value = 8497.5 # K126: was 8568.856800000001
$ws.Cells.Item(126, 12).Value = 12298.9995 # L126: was 12213
$ws.Cells.Item(126, 13).Value = -6027.5 # M126: was -6098.856800000001
$ws.Cells.Item(126, 14).Value = -17238.9995 # N126: was -17153

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3842.0588 # H40: was 3607.95
$ws.Cells.Item(40, 9).Value = 3025.4614 # I40: was 2885.9375
$ws.Cells.Item(40, 11).Value = 3025.4614 # K40: was 2885.9375
$ws.Cells.Item(40, 13).Value = -2889.4614 # M40: was -2749.9375

$ws.Cells.Item(46, 8).Value = 41428.816 # H46: was 30807.666
$ws.Cells.Item(46, 9).Value = 106979.75 # I46: was 85783.8
$ws.Cells.Item(46, 10).Value = 3971.1428 # J46: was 3319.6
$ws.Cells.Item(46, 11).Value = 106979.75 # K46: was 85783.8
$ws.Cells.Item(46, 12).Value = 3971.1428 # L46: was 3319.6
$ws.Cells.Item(46, 13).Value = -106791.75 # M46: was -85595.8
$ws.Cells.Item(46, 14).Value = -4347.1428 # N46: was -3695.6

$ws.Cells.Item(61, 8).Value = 396308 # H61: was 527777.7
$ws.Cells.Item(61, 9).Value = 417299.66 # I61: was 625000
$ws.Cells.Item(61, 11).Value = 417299.66 # K61: was 625000
$ws.Cells.Item(61, 13).Value = -417097.66 # M61: was -624798

$ws.Cells.Item(113, 8).Value = 396308 # H113: was 527777.7
$ws.Cells.Item(113, 9).Value = 417299.66 # I113: was 625000
$ws.Cells.Item(113, 11).Value = 417299.66 # K113: was 625000
$ws.Cells.Item(113, 13).Value = -415129.66 # M113: was -622830

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(105, 8).Value = 40615 # H105: was 40607.5
$ws.Cells.Item(105, 10).Value = 40615 # J105: was 40607.5
$ws.Cells.Item(105, 12).Value = 40615 # L105: was 40607.5
$ws.Cells.Item(105, 14).Value = -47603 # N105: was -47595.5

$ws.Cells.Item(122, 8).Value = 1283.7778 # H122: was 1360.9546
$ws.Cells.Item(122, 9).Value = 1283.7778 # I122: was 1306.3529
$ws.Cells.Item(122, 10).Value = 0 # J122: was 1546.6
$ws.Cells.Item(122, 11).Value = 3851.3334 # K122: was 3919.0587
$ws.Cells.Item(122, 12).Value = 0 # L122: was 4639.799999999999
$ws.Cells.Item(122, 13).Value = -1401.3334 # M122: was -1469.0587
$ws.Cells.Item(122, 14).ClearContents() # N122: was -9539.799999999999

$ws.Cells.Item(132, 8).Value = 1071810.1 # H132: was 1132956.5
$ws.Cells.Item(132, 9).Value = 1278762.4 # I132: was 1321361.1
$ws.Cells.Item(132, 10).Value = 2557.1667 # J132: was 2528.8
$ws.Cells.Item(132, 11).Value = 3836287.2 # K132: was 3964083.3
$ws.Cells.Item(132, 12).Value = 7671.500100000001 # L132: was 7586.400000000001
$ws.Cells.Item(132, 13).Value = -3833757.2 # M132: was -3961553.3
$ws.Cells.Item(132, 14).Value = -12731.5001 # N132: was -12646.4
